$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect (with the known password) so the
# holdings data and the confidential-notice text can be updated, then
# re-apply protection with the same password when done.
$ws.Unprotect("D382")

# --- Update the "as of" date in the confidential notice (A41) ----------
$notice = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + "`n" + "Model holdings provided as of 2021-05-19 for illustrative purposes only and are subject to change."
$ws.Range("A41").Value = $notice

# --- Refresh the Weight (D) / Percent Change (E) columns ----------------
$ws.Range("D2").Value = 0.03248966668065826
$ws.Range("E2").Value = -0.01160310553707022
$ws.Range("D3").Value = 0.0285479973674686
$ws.Range("E3").Value = -0.0003120968748699937
$ws.Range("D4").Value = 0.02800351320977965
$ws.Range("E4").Value = 0.02404621171413224
$ws.Range("D5").Value = 0.06399728193508482
$ws.Range("E5").Value = -0.000148501986214078
$ws.Range("D6").Value = 0.01579281249232772
$ws.Range("E6").Value = -0.002707988566270414
$ws.Range("D7").Value = 0.01531406242203977
$ws.Range("E7").Value = 0.01357536265611659
$ws.Range("D8").Value = 0.02980130090134898
$ws.Range("E8").Value = -0.008889420393845215
$ws.Range("D9").Value = 0.0348826250551414
$ws.Range("E9").Value = -0.006016573958451676
$ws.Range("D10").Value = 0.02888419158556163
$ws.Range("E10").Value = 0.002920128321131754
$ws.Range("D11").Value = 0.02962073015523541
$ws.Range("E11").Value = -0.005521242747520083
$ws.Range("D12").Value = 0.01093225191519828
$ws.Range("E12").Value = -0.00642941229738303
$ws.Range("D13").Value = 0.01442229636886485
$ws.Range("E13").Value = -0.002498558523928551
$ws.Range("D14").Value = 0.01422410413546608
$ws.Range("E14").Value = 0.001266686154146157
$ws.Range("D15").Value = 0.009125950471345087
$ws.Range("E15").Value = -0.004859845526338669
$ws.Range("D16").Value = 0.008170232278921614
$ws.Range("E16").Value = 0.0004362050163575493
$ws.Range("D17").Value = 0.02970527369535656
$ws.Range("E17").Value = 0.03270657397471188
$ws.Range("D18").Value = 0.02522050618415207
$ws.Range("E18").Value = -0.003540587219343694
$ws.Range("D19").Value = 0.0334828057842829
$ws.Range("E19").Value = 0.001773993258825657
$ws.Range("D20").Value = 0.03068514718495751
$ws.Range("E20").Value = 0.01171118854045683
$ws.Range("D21").Value = 0.04560658703074375
$ws.Range("E21").Value = 0.002292233755746942
$ws.Range("D22").Value = 0.03609038991401507
$ws.Range("E22").Value = -0.01695194206714934
$ws.Range("D23").Value = 0.03135733762690441
$ws.Range("E23").Value = -0.00726124704025255
$ws.Range("D24").Value = 0.03090432680770721
$ws.Range("E24").Value = -0.01239693247996321
$ws.Range("D25").Value = 0.01514299539940586
$ws.Range("E25").Value = -0.01164979995293025
$ws.Range("D26").Value = 0.01498499599655649
$ws.Range("E26").Value = 0.001268431901062295
$ws.Range("D27").Value = 0.0313426860532066
$ws.Range("E27").Value = -0.00500312695434646
$ws.Range("D28").Value = 0.03124982675504074
$ws.Range("E28").Value = 0.007349586902529293
$ws.Range("D29").Value = 0.02887706379295188
$ws.Range("E29").Value = 0.0001645548790520923
$ws.Range("D30").Value = 0.02952410896652551
$ws.Range("E30").Value = -0.01917969902626138
$ws.Range("D31").Value = 0.03330045309001689
$ws.Range("E31").Value = 0.003567415229295712
$ws.Range("D32").Value = 0.03169115591412752
$ws.Range("E32").Value = -0.0005497903924127634
$ws.Range("D33").Value = 0.02852938590898759
$ws.Range("E33").Value = 0.0140743413929989
$ws.Range("D34").Value = 0.03250115034652952
$ws.Range("E34").Value = -0.002144354013353444
$ws.Range("D35").Value = 0.0307528612147501
$ws.Range("E35").Value = -0.0003476648510835068
$ws.Range("D36").Value = 0.03214436472756389
$ws.Range("E36").Value = -0.007637819525715961
$ws.Range("D37").Value = 0.03269756063177586
$ws.Range("E37").Value = 0.06094075473526117
$ws.Range("E38").Value = 0.001466543329455394

# Re-protect the worksheet with the original password.
$ws.Protect("D382", $true, $true, $true)
